$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (row 1) to match the new summer-reporting template
$ws.Range("E1").Value = "Legal Given Name"
$ws.Range("F1").Value = "Birthdate"
$ws.Range("G1").Value = "Ministry Course Code and Level"
$ws.Range("I1").Value = "Final Percent"
$ws.Range("K1").Value = "Credits"

# Update the selected range shown when the sheet is opened
$ws.Range("A1:K1").Select()
